$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells that look numeric as Text so Excel keeps them as strings
$priceTextCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.338.95"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "1.936.76"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "0.7667"
$ws.Range("E5").Value = "  +5.94%  "

$ws.Range("D6").Value = "248.57"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").Value = "0.9974"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +1.31%  "

$ws.Range("D9").Value = "0.3232"
$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("D10").Value = "0.07145"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "0.7933"
$ws.Range("E11").Value = "  -2.26%  "

$ws.Range("D12").Value = "0.08004"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").Value = "1.931.10"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "5.393"
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "94.88"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "14.78"
$ws.Range("E16").Value = "  -2.76%  "

$ws.Range("D17").Value = "30.332.49"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "254.08"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "0.000008058"
$ws.Range("E19").Value = "  -3.01%  "

$ws.Range("D20").Value = "5.814"
$ws.Range("E20").Value = "  -1.58%  "

$ws.Range("D21").Value = "2.188.37"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "0.9968"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").Value = "0.9968"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").Value = "6.848"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").Value = "9.614"
$ws.Range("E25").Value = "  -1.41%  "

$ws.Range("D26").Value = "165.35"
$ws.Range("E26").Value = "  +1.07%  "

$ws.Range("D27").Value = "0.1362"
$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").Value = "2.327"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("D29").Value = "19.14"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("D32").Value = "4.444"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "4.158"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").Value = "0.05198"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").Value = "1.298"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").Value = "0.7560"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("D37").Value = "2.767"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "0.01976"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("D39").Value = "2.801"
$ws.Range("E39").Value = "  -1.16%  "

$ws.Range("D40").Value = "78.60"
$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("D41").Value = "6.450"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("D42").Value = "0.4542"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "2.003"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("D44").Value = "0.9982"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").Value = "0.8372"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D46").Value = "102.62"
$ws.Range("E46").Value = "  +0.69%  "

$ws.Range("D47").Value = "7.576"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("D48").Value = "9.843"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "990.14"
$ws.Range("E49").Value = "  +13.26%  "

$ws.Range("D50").Value = "37.47"
$ws.Range("E50").Value = "  +1.66%  "

$ws.Range("D51").Value = "0.4197"
$ws.Range("E51").Value = "  +0.32%  "
